$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.643.25"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.05%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.299.68"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.09%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("B5").Value = "BNB"
$ws.Range("C5").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "268.02"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.22%  "
$ws.Range("B6").Value = "Solana"
$ws.Range("C6").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "94.54"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +8.35%  "
$ws.Range("E7").Value = "  +0.38%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.622"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.28%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "44.65"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.56%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0937"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.01%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.08"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +6.29%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.104"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.59%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.643.54"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.94%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.35"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.01%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.860"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +7.84%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.301.69"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.71%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.590.99"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.99%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0000108"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.33%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.33"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.24%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "71.19"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.96%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.29"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.08%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "237.75"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.59%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.56"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +7.57%  "
$ws.Range("E25").Value = "  -0.01%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.33"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.21%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.50"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.30%  "
$ws.Range("B28").Value = "WEMIXToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.39"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.21%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.30"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.01%  "
$ws.Range("E30").Value = "  -3.06%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "22.29"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +6.70%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "171.71"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.00%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0895"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.22%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.55"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.20%  "
$ws.Range("E35").Value = "  +1.72%  "
$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.109"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.44%  "
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0356"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.80%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.43"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.49%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.45"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.04%  "
$ws.Range("B40").Value = "LidoDAOToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.31"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.03%  "
$ws.Range("B41").Value = "Algorand"
$ws.Range("C41").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.233"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +14.69%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.37"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +19.84%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "12.03"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.95%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.44"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.05%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "61.83"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.69%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.08"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +6.01%  "
$ws.Range("E47").Value = "  +3.52%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "100.47"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.63%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.21"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.50%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.523.45"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.93%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.421"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.31%  "
